$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.975738763809204
$ws.Range("B1").Value = 2.697833776473999
$ws.Range("C1").Value = 2.135915040969849
$ws.Range("D1").Value = 2.004388809204102
$ws.Range("E1").Value = 1.747995257377625
